$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain (unambiguous) text -- assign directly.
$ws.Range("D2").Value = '28.616.28'
$ws.Range("E2").Value = '  +2.50%  '
$ws.Range("D3").Value = '1.913.01'
$ws.Range("E3").Value = '  +5.60%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("E10").Value = '  +5.44%  '
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("E13").Value = '  +3.18%  '
$ws.Range("D14").Value = '1.916.40'
$ws.Range("E14").Value = '  +5.66%  '
$ws.Range("E15").Value = '  +3.52%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("E20").Value = '  +5.66%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("E22").Value = '  +5.86%  '
$ws.Range("D23").Value = '28.660.33'
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("E26").Value = '  +15.47%  '
$ws.Range("D27").Value = '2.134.51'
$ws.Range("E27").Value = '  +5.56%  '
$ws.Range("E28").Value = '  +3.90%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  +7.41%  '
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("E33").Value = '  +2.74%  '
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("E35").Value = '  +9.21%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("E37").Value = '  +4.97%  '
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("E39").Value = '  +3.56%  '
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("E42").Value = '  +4.85%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  +4.97%  '
$ws.Range("E45").Value = '  +3.97%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("E48").Value = '  +5.79%  '
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("E50").Value = '  +2.99%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E51").Value = '  +6.57%  '

# Cells whose new values look like numbers (e.g. "1.001") but must stay as
# plain text, matching the rest of the sheet -- force Text format first so
# Excel does not reinterpret the string as a numeric value.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.70'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5046'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3966'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09632'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.161'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.62'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.567'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.18'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.556'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001136'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.91'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06636'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.06'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9996'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.279'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.44'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.287'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.763'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.41'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.36'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.89'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.117'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1077'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.727'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.633'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.821'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06785'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02442'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2224'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.110'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.64'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6426'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.198'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.80'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6107'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.287'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.666'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.044'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.16'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.210'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.37'
